$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row renames (row 1) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case 'de/del/la/las/el/los/y' connector words in state/municipality names ---
$ws.Range('B7').Value = 'Pabellón De Arteaga'
$ws.Range('B8').Value = 'Rincón De Romos'
$ws.Range('B9').Value = 'San José De Gracia'
$ws.Range('B14').Value = 'Playas De Rosarito'
$ws.Range('B74').Value = 'Coyame Del Sotol'
$ws.Range('B85').Value = 'Guadalupe Y Calvo'
$ws.Range('B88').Value = 'Hidalgo Del Parral'
$ws.Range('B113').Value = 'San Francisco De Borja'
$ws.Range('B114').Value = 'San Francisco De Conchos'
$ws.Range('B115').Value = 'San Francisco Del Oro'
$ws.Range('B123').Value = 'Valle De Zaragoza'
$ws.Range('B141').Value = 'San Juan De Sabinas'
$ws.Range('A154').Value = 'Ciudad De México'
$ws.Range('B168').Value = 'Coneto De Comonfort'
$ws.Range('B182').Value = 'Nombre De Dios'
$ws.Range('B186').Value = 'Pánuco De Coronado'
$ws.Range('B193').Value = 'San Juan De Guadalupe'
$ws.Range('B194').Value = 'San Juan Del Río'
$ws.Range('B195').Value = 'San Luis Del Cordero'
$ws.Range('B196').Value = 'San Pedro Del Gallo'
$ws.Range('A206').Value = 'Estado De México'
$ws.Range('B206').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B209').Value = 'Almoloya De Alquisiras'
$ws.Range('B210').Value = 'Almoloya De Juárez'
$ws.Range('B222').Value = 'Ecatepec De Morelos'
$ws.Range('B223').Value = 'Ixtapan De La Sal'
$ws.Range('B230').Value = 'Naucalpan De Juárez'
$ws.Range('B234').Value = 'San Felipe Del Progreso'
$ws.Range('B235').Value = 'San Martín De Las Pirámides'
$ws.Range('B243').Value = 'Tenango Del Valle'
$ws.Range('B249').Value = 'Tlalnepantla De Baz'
$ws.Range('B254').Value = 'Valle De Bravo'
$ws.Range('B255').Value = 'Villa De Allende'
$ws.Range('B266').Value = 'San Miguel De Allende'
$ws.Range('B267').Value = 'Apaseo El Alto'
$ws.Range('B268').Value = 'Apaseo El Grande'
$ws.Range('B273').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B277').Value = 'Jaral Del Progreso'
$ws.Range('B282').Value = 'Purísima Del Rincón'
$ws.Range('B286').Value = 'San Diego De La Unión'
$ws.Range('B288').Value = 'San Francisco Del Rincón'
$ws.Range('B290').Value = 'San Luis De La Paz'
$ws.Range('B291').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B292').Value = 'Silao De La Victoria'
$ws.Range('B295').Value = 'Valle De Santiago'
$ws.Range('B301').Value = 'Acapulco De Juárez'
$ws.Range('B303').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B306').Value = 'Atenango Del Río'
$ws.Range('B307').Value = 'Atoyac De Álvarez'
$ws.Range('B308').Value = 'Ayutla De Los Libres'
$ws.Range('B310').Value = 'Chilapa De Álvarez'
$ws.Range('B311').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B312').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B315').Value = 'Coyuca De Benítez'
$ws.Range('B316').Value = 'Coyuca De Catalán'
$ws.Range('B318').Value = 'Cutzamala De Pinzón'
$ws.Range('B322').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B323').Value = 'Iguala De La Independencia'
$ws.Range('B324').Value = 'Zihuatanejo De Azueta'
$ws.Range('B333').Value = 'Taxco De Alarcón'
$ws.Range('B335').Value = 'Técpan De Galeana'
$ws.Range('B337').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B340').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B341').Value = 'Tlapa De Comonfort'
$ws.Range('B357').Value = 'Huejutla De Reyes'
$ws.Range('B360').Value = 'Jacala De Ledezma'
$ws.Range('B362').Value = 'Mixquiahuala De Juárez'
$ws.Range('B363').Value = 'Nopala De Villagrán'
$ws.Range('B364').Value = 'Pachuca De Soto'
$ws.Range('B365').Value = 'Progreso De Obregón'
$ws.Range('B367').Value = 'Santiago De Anaya'
$ws.Range('B368').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B373').Value = 'Tula De Allende'
$ws.Range('B374').Value = 'Tulancingo De Bravo'
$ws.Range('B379').Value = 'Acatlán De Juárez'
$ws.Range('B380').Value = 'Ahualulco De Mercado'
$ws.Range('B385').Value = 'Atotonilco El Alto'
$ws.Range('B386').Value = 'Autlán De Navarro'
$ws.Range('B395').Value = 'Concepción De Buenos Aires'
$ws.Range('B399').Value = 'Encarnación De Díaz'
$ws.Range('B404').Value = 'Huejuquilla El Alto'
$ws.Range('B405').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B406').Value = 'Ixtlahuacán Del Río'
$ws.Range('B413').Value = 'Lagos De Moreno'
$ws.Range('B418').Value = 'Ojuelos De Jalisco'
$ws.Range('B420').Value = 'San Diego De Alejandría'
$ws.Range('B424').Value = 'San Miguel El Alto'
$ws.Range('B425').Value = 'Santa María De Los Ángeles'
$ws.Range('B428').Value = 'Tamazula De Gordiano'
$ws.Range('B432').Value = 'Teocuitatlán De Corona'
$ws.Range('B433').Value = 'Tepatitlán De Morelos'
$ws.Range('B435').Value = 'Tizapán El Alto'
$ws.Range('B436').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B444').Value = 'Unión De Tula'
$ws.Range('B445').Value = 'Valle De Guadalupe'
$ws.Range('B450').Value = 'Yahualica De González Gallo'
$ws.Range('B451').Value = 'Zacoalco De Torres'
$ws.Range('B454').Value = 'Zapotlán Del Rey'
$ws.Range('B455').Value = 'Zapotlán El Grande'
$ws.Range('B469').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B471').Value = 'Cojumatlán De Régules'
$ws.Range('B509').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B533').Value = 'Puente De Ixtla'
$ws.Range('B535').Value = 'Tetela Del Volcán'
$ws.Range('B536').Value = 'Tlaltizapán De Zapata'
$ws.Range('B543').Value = 'Amatlán De Cañas'
$ws.Range('B544').Value = 'Bahía De Banderas'
$ws.Range('B548').Value = 'Ixtlán Del Río'
$ws.Range('B552').Value = 'Santa María Del Oro'
$ws.Range('B568').Value = 'Mier Y Noriega'
$ws.Range('B570').Value = 'San Nicolás De Los Garza'
$ws.Range('B572').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B573').Value = 'Cuilápam De Guerrero'
$ws.Range('B575').Value = 'Guevea De Humboldt'
$ws.Range('B576').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B577').Value = 'Ixtlán De Juárez'
$ws.Range('B578').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B580').Value = 'Mariscala De Juárez'
$ws.Range('B582').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B583').Value = 'Nejapa De Madero'
$ws.Range('B584').Value = 'Oaxaca De Juárez'
$ws.Range('B585').Value = 'Ocotlán De Morelos'
$ws.Range('B621').Value = 'Tepelmeme Villa De Morelos'
$ws.Range('B622').Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range('B623').Value = 'Tlacolula De Matamoros'
$ws.Range('B634').Value = 'Chalchicomula De Sesma'
$ws.Range('B642').Value = 'Huehuetlán El Grande'
$ws.Range('B645').Value = 'Ixcamilpa De Guerrero'
$ws.Range('B647').Value = 'Izúcar De Matamoros'
$ws.Range('B652').Value = 'Los Reyes De Juárez'
$ws.Range('B653').Value = 'Palmar De Bravo'
$ws.Range('B660').Value = 'San Nicolás De Los Ranchos'
$ws.Range('B663').Value = 'San Salvador El Seco'
$ws.Range('B664').Value = 'San Salvador El Verde'
$ws.Range('B670').Value = 'Tepexi De Rodríguez'
$ws.Range('B672').Value = 'Tetela De Ocampo'
$ws.Range('B673').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B685').Value = 'Amealco De Bonfil'
$ws.Range('B687').Value = 'Cadereyta De Montes'
$ws.Range('B691').Value = 'Jalpan De Serra'
$ws.Range('B692').Value = 'Landa De Matamoros'
$ws.Range('B693').Value = 'Pinal De Amoles'
$ws.Range('B696').Value = 'San Juan Del Río'
$ws.Range('B707').Value = 'Ciudad Del Maíz'
$ws.Range('B716').Value = 'San Ciro De Acosta'
$ws.Range('B718').Value = 'Santa María Del Río'
$ws.Range('B723').Value = 'Villa De Arista'
$ws.Range('B724').Value = 'Villa De Guadalupe'
$ws.Range('B725').Value = 'Villa De Ramos'
$ws.Range('B726').Value = 'Villa De Reyes'
$ws.Range('B768').Value = 'Nacozari De García'
$ws.Range('B776').Value = 'San Pedro De La Cueva'
$ws.Range('B810').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('B825').Value = 'Amatlán De Los Reyes'
$ws.Range('B835').Value = 'Cosamaloapan De Carpio'
$ws.Range('B842').Value = 'Hueyapan De Ocampo'
$ws.Range('B852').Value = 'Lerdo De Tejada'
$ws.Range('B853').Value = 'Martínez De La Torre'
$ws.Range('B862').Value = 'Poza Rica De Hidalgo'
$ws.Range('B869').Value = 'Sayula De Alemán'
$ws.Range('B878').Value = 'Vega De Alatorre'
$ws.Range('B890').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B892').Value = 'Concepción Del Oro'
$ws.Range('B910').Value = 'Moyahua De Estrada'
$ws.Range('B911').Value = 'Nochistlán De Mejía'
$ws.Range('B921').Value = 'Teúl De González Ortega'
$ws.Range('B922').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B925').Value = 'Villa De Cos'

# --- Minor recompute of D625 percentage (Total row for Oaxaca sub-group) ---
$ws.Range('D625').Value = 0.009792284866468845

# --- Remove trailing footnote/metadata rows 934-938 (dimension shrinks to A1:D932) ---
$ws.Rows("934:938").Delete()

Write-Output 'done'